$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update PO values for the Farm row (row 3)
$ws.Range("D3").Value = -10
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 40

# Update the active selection to F3 as shown in the diff
$ws.Range("F3").Select()
